# Update average_county_temperature (column I) with NOAA-sourced data,
# and the dependent worst/best ASHP COP values (columns N/O) for the
# affected rows (facilities 1002036, 1002674, 1011252, 1013683).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("I3").Value = -1.819444444444444
$ws.Range("N3").Value = 1.181188907392658
$ws.Range("O3").Value = 1.232403941639294

# Row 4
$ws.Range("I4").Value = 0.2777777777777778
$ws.Range("N4").Value = 1.194245973645681
$ws.Range("O4").Value = 1.246785162287481

# Row 8
$ws.Range("I8").Value = 19.79629629629628
$ws.Range("N8").Value = 1.331198999020781
$ws.Range("O8").Value = 1.39868801294648

# Row 9
$ws.Range("I9").Value = 21.28240740740739
$ws.Range("N9").Value = 1.342924567132234
$ws.Range("O9").Value = 1.411784266254412
